# Regenerate the "K" column (column G) values for rows 2-22 on Sheet1.
# These are the recalculated strikeout (K) counts that replace the old
# "Strike#" derived values, per the commit message:
#   "regen save_data to use K instead of Strike#, regen std/mean,
#    calc and write s_vals"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @(0, 0, 1, 1, 1, 0, 2, 0, 1, 1, 2, 2, 3, 2, 0, 0, 1, 1, 1, 3, 1)

for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
